$d = $word.ActiveDocument

# Locate the existing closing paragraph ("Please consult the rubric ...")
# so the edit is anchored to stable content rather than a fragile index.
$anchorText = "Please consult the rubric both prior to and throughout the report-writing process to ensure your work aligns with the stated evaluation criteria."
$findRange = $d.Content
$found = $findRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph text."
}

# Expand to the whole paragraph, including its trailing paragraph mark, so
# the replacement XML below fully owns (and re-creates) that paragraph mark
# instead of leaving a stray empty leftover paragraph behind.
$para = $findRange.Paragraphs(1)
$targetRange = $para.Range

# Replace that paragraph with itself (now carrying paragraph-level Arial/12pt
# formatting on the mark, matching a paragraph that was extended by pressing
# Enter) followed by the new paragraphs that were appended to the document:
# a couple of blank spacer paragraphs, one with a double bottom border acting
# as a divider, and then a short pasted-in Q&A conversation about Assignment 2.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Please consult the rubric both prior to and throughout the report-writing process to ensure your work aligns with the stated evaluation criteria.</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>For assignment 2, choosing option 1, all you want us to do is get a screenrecording of us visiting shoden and finding open cctv?</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Display the search used, showcase the camera including the view from the camera, and any details you can (for example camera hardware and software, version, geolocation data, etc.)</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve">More on this ^, How long do you want the recording to last? </w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>and is there an example video to watch; to see what you are looking for?</w:t></w:r></w:p>
<w:p><w:r><w:t>As long as needed to show the details or the hack</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

# Do not include the final paragraph mark in the replaced range: Word always
# needs a terminating paragraph mark for the document/story, so leaving the
# original one in place (and letting InsertXML's content flow into it) avoids
# an extra blank paragraph being left behind after the inserted content.
$targetRange.End = $targetRange.End - 1
$targetRange.InsertXML($xml)
